$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has a "2019" data column in P (years 2007..2019 laid out
# across D:P, one column per year). We are adding a new "2020" column in Q,
# re-using the same per-row formatting as column P.

# 1) Copy the formatting of the whole P3:P33 block into Q3:Q33 in a single
#    operation (doing this range-wise, rather than cell-by-cell, keeps each
#    destination cell's own style instead of falling back to the row's
#    default style).
$ws.Range("P3:P33").Copy() | Out-Null
$ws.Range("Q3:Q33").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# 2) Header: year 2020
$ws.Cells.Item(3, 17).Value = 2020

# 3) Data values for 2020, row by row (same ordering as column P).
$ws.Cells.Item(4, 17).Value = 1.9148453093736542
$ws.Cells.Item(5, 17).Value = 1.7453236044300597
$ws.Cells.Item(6, 17).Value = 2.0818900906859255
$ws.Cells.Item(7, 17).Value = 1.658050942694075
$ws.Cells.Item(8, 17).Value = 1.4467487937731931
$ws.Cells.Item(9, 17).Value = 1.8774124750304142
$ws.Cells.Item(10, 17).Value = 0.96024351775610284
$ws.Cells.Item(11, 17).Value = 0.63595936855594293
$ws.Cells.Item(12, 17).Value = 1.2888424905592288
$ws.Cells.Item(13, 17).Value = 1.6032353288937073
$ws.Cells.Item(14, 17).Value = 2.4146715443031859
$ws.Cells.Item(15, 17).Value = 0.79837132250209564
$ws.Cells.Item(16, 17).Value = 1.3751327862596732
$ws.Cells.Item(17, 17).Value = 0.67516929870164943
$ws.Cells.Item(18, 17).Value = 2.1012817818869509
$ws.Cells.Item(19, 17).Value = 1.5943738893736428
$ws.Cells.Item(20, 17).Value = 1.5765365498500856
$ws.Cells.Item(21, 17).Value = 1.6126194804433236
$ws.Cells.Item(22, 17).Value = 0.37150276583809166
$ws.Cells.Item(23, 17).Value = 0
$ws.Cells.Item(24, 17).Value = 0.75125835774923
$ws.Cells.Item(25, 17).Value = 2.8942542850468351
$ws.Cells.Item(26, 17).Value = 2.72898263527357
$ws.Cells.Item(27, 17).Value = 3.0545792215303034
$ws.Cells.Item(28, 17).Value = 3.9473869708034344
$ws.Cells.Item(29, 17).Value = 3.6031203021816895
$ws.Cells.Item(30, 17).Value = 4.2520923837938582
$ws.Cells.Item(31, 17).Value = 0
$ws.Cells.Item(32, 17).Value = 0
$ws.Cells.Item(33, 17).Value = 0

# 4) Match the cursor location left by the author in the saved file.
$ws.Range("T1").Select() | Out-Null
